# lambrechts 6653: use cellranger instead of dropSeqPipe.
#
# The sample-name extraction formula in column A used to split on the first
# underscore ("_"), which worked for names like "scrBT1429m_S0_..." but broke
# for the cellranger-produced names like "scrBT1425_hg19_S11_..." (it cut the
# name down to "scrBT1425" instead of "scrBT1425_hg19"). Update the formula
# to split on "_S" (the start of the Sample index token) instead, and refresh
# the dependent values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-fill the shared formula for A2:A12 with the corrected split token, then
# set A13 (which keeps the corrected formula too but no longer shares the
# same formula group, matching the prior row-13 one-off edit).
$ws.Range("A2:A12").Formula = '=LEFT(AT2,FIND("_S",AT2)-1)'
$ws.Range("A13").Formula = '=LEFT(AT13,FIND("_S",AT13)-1)'

# Column A now holds longer sample names (e.g. "scrBT1425_hg19"); widen it
# and drop the stale auto "best fit" flag.
$ws.Columns.Item(1).ColumnWidth = 25.43

# Restore the cursor/selection position left after the edit.
$ws.Range("B17").Select()
